# Case and Fatality Demographics Data Updated
# Refresh the three "Fatalities by ..." breakdown tabs with the latest
# counts (percent columns are formulas and recalculate automatically),
# and restore each sheet's last-used selection, moving the selected tab
# from "Fatalities by Age Group" to "Fatalities by Race-Ethnicity".

$wb = $excel.ActiveWorkbook

# --- Fatalities by Age Group ---
$wsFA = $wb.Worksheets.Item("Fatalities by Age Group")
$wsFA.Range("B3").Value = 14
$wsFA.Range("B4").Value = 36
$wsFA.Range("B5").Value = 280
$wsFA.Range("B6").Value = 925
$wsFA.Range("B7").Value = 2689
$wsFA.Range("B8").Value = 5979
$wsFA.Range("B9").Value = 4940
$wsFA.Range("B10").Value = 6332
$wsFA.Range("B11").Value = 6949
$wsFA.Range("B12").Value = 6837
$wsFA.Range("B13").Value = 17053

# --- Fatalities by Gender ---
$wsFG = $wb.Worksheets.Item("Fatalities by Gender")
$wsFG.Range("B2").Value = 21826
$wsFG.Range("B3").Value = 30214

# --- Fatalities by Race-Ethnicity ---
$wsFR = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$wsFR.Range("B2").Value = 1099
$wsFR.Range("B3").Value = 5319
$wsFR.Range("B4").Value = 24101
$wsFR.Range("B5").Value = 287
$wsFR.Range("B6").Value = 21212

# --- Selections on each sheet (tracked per-sheet, last one set becomes the
#     active cell shown when that sheet is later activated) ---
$wsCA = $wb.Worksheets.Item("Cases by Age Group")
$wsCA.Range("E14").Select() | Out-Null

$wsCG = $wb.Worksheets.Item("Cases by Gender")
$wsCG.Range("B2:B4").Select() | Out-Null

$wsCR = $wb.Worksheets.Item("Cases by RaceEthnicity")
$wsCR.Range("I10").Select() | Out-Null

$wsFA.Range("H10").Select() | Out-Null

$wsFG.Range("F15").Select() | Out-Null

$wsFR.Range("D10").Select() | Out-Null

# --- Active sheet / tab selection moves from "Fatalities by Age Group" to
#     "Fatalities by Race-Ethnicity" ---
$wsFR.Activate() | Out-Null
